$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Table1")
$col = $tbl.ListColumns.Add()
$ws.Cells.Item(1,10).Value = "Authorship Resource"
$ws.Range("J2:J33").Value = "Noémi Villars-Amberg, Daniela Subotic"
